$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers
$ws.Range("A1").Value = "Nom et prénom"
$ws.Range("B1").Value = "grammaire"
$ws.Range("C1").Value = "expression de text"
$ws.Range("D1").Value = "resume"
$ws.Range("E1").Value = "fanadihadiana lahatsoratra"
$ws.Range("F1").Value = "famakafakan-kevitra"
$ws.Range("G1").Value = "probabilite"
$ws.Range("H1").Value = "statistique"
$ws.Range("I1").Value = "geographie"
$ws.Range("J1").Value = "histoire"
$ws.Range("K1").Value = "geologie"
$ws.Range("L1").Value = "suite numerique"
$ws.Range("M1").Value = "fonction"
$ws.Range("N1").Value = "interference mecanique"
$ws.Range("O1").Value = "reproduction humaine"
$ws.Range("P1").Value = "vibration sonore"
$ws.Range("Q1").Value = "interference lumineuse"
$ws.Range("R1").Value = "genetique"

# Row 2 values
$ws.Range("A2").Value = "RAMANANDRAIBE"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "]15,20]"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "[10]"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = "]15,20]"
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "]15,20]"
